$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Mdk"
$ws.Cells.Item(2,3).Value = "Itga4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.098888
$ws.Cells.Item(2,8).Value = 6.296664
$ws.Cells.Item(2,9).Value = 0.1082453658858517
$ws.Cells.Item(2,10).Value = 0.1082453658858517
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 26.04517333333333
$ws.Cells.Item(2,14).Value = 78.13552
$ws.Cells.Item(2,15).Value = 0.9210237118384171
$ws.Cells.Item(2,16).Value = 0.921023711838417
$ws.Cells.Item(2,17).Value = 54.66590176725334
$ws.Cells.Item(2,18).Value = 491.99311590528
$ws.Cells.Item(2,19).Value = 0.09969654867749468
$ws.Cells.Item(2,20).Value = 0.09969654867749465
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Mdk"
$ws.Cells.Item(3,3).Value = "Itga4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.098888
$ws.Cells.Item(3,8).Value = 6.296664
$ws.Cells.Item(3,9).Value = 0.1082453658858517
$ws.Cells.Item(3,10).Value = 0.1082453658858517
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.3302223333333333
$ws.Cells.Item(3,14).Value = 0.990667
$ws.Cells.Item(3,15).Value = 0.01167750336256582
$ws.Cells.Item(3,16).Value = 0.01167750336256582
$ws.Cells.Item(3,17).Value = 0.6930996927653333
$ws.Cells.Item(3,18).Value = 6.237897234888
$ws.Cells.Item(3,19).Value = 0.001264035624114201
$ws.Cells.Item(3,20).Value = 0.001264035624114201
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Mdk"
$ws.Cells.Item(4,3).Value = "Itga4"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.098888
$ws.Cells.Item(4,8).Value = 6.296664
$ws.Cells.Item(4,9).Value = 0.1082453658858517
$ws.Cells.Item(4,10).Value = 0.1082453658858517
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.903109
$ws.Cells.Item(4,14).Value = 5.709327
$ws.Cells.Item(4,15).Value = 0.06729878479901708
$ws.Cells.Item(4,16).Value = 0.06729878479901708
$ws.Cells.Item(4,17).Value = 3.994412642792
$ws.Cells.Item(4,18).Value = 35.949713785128
$ws.Cells.Item(4,19).Value = 0.007284781584242796
$ws.Cells.Item(4,20).Value = 0.007284781584242794
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Mdk"
$ws.Cells.Item(5,3).Value = "Itga4"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 15.87514366666667
$ws.Cells.Item(5,8).Value = 47.625431
$ws.Cells.Item(5,9).Value = 0.8187243600843848
$ws.Cells.Item(5,10).Value = 0.8187243600843847
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 26.04517333333333
$ws.Cells.Item(5,14).Value = 78.13552
$ws.Cells.Item(5,15).Value = 0.9210237118384171
$ws.Cells.Item(5,16).Value = 0.921023711838417
$ws.Cells.Item(5,17).Value = 413.4708684899022
$ws.Cells.Item(5,18).Value = 3721.23781640912
$ws.Cells.Item(5,19).Value = 0.7540645490974529
$ws.Cells.Item(5,20).Value = 0.7540645490974527
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Mdk"
$ws.Cells.Item(6,3).Value = "Itga4"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 15.87514366666667
$ws.Cells.Item(6,8).Value = 47.625431
$ws.Cells.Item(6,9).Value = 0.8187243600843848
$ws.Cells.Item(6,10).Value = 0.8187243600843847
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.3302223333333333
$ws.Cells.Item(6,14).Value = 0.990667
$ws.Cells.Item(6,15).Value = 0.01167750336256582
$ws.Cells.Item(6,16).Value = 0.01167750336256582
$ws.Cells.Item(6,17).Value = 5.242326983608556
$ws.Cells.Item(6,18).Value = 47.180942852477
$ws.Cells.Item(6,19).Value = 0.009560656467899956
$ws.Cells.Item(6,20).Value = 0.009560656467899954
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Mdk"
$ws.Cells.Item(7,3).Value = "Itga4"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 15.87514366666667
$ws.Cells.Item(7,8).Value = 47.625431
$ws.Cells.Item(7,9).Value = 0.8187243600843848
$ws.Cells.Item(7,10).Value = 0.8187243600843847
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.903109
$ws.Cells.Item(7,14).Value = 5.709327
$ws.Cells.Item(7,15).Value = 0.06729878479901708
$ws.Cells.Item(7,16).Value = 0.06729878479901708
$ws.Cells.Item(7,17).Value = 30.21212878832633
$ws.Cells.Item(7,18).Value = 271.909159094937
$ws.Cells.Item(7,19).Value = 0.05509915451903197
$ws.Cells.Item(7,20).Value = 0.05509915451903197
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Mdk"
$ws.Cells.Item(8,3).Value = "Itga4"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.416064
$ws.Cells.Item(8,8).Value = 4.248192
$ws.Cells.Item(8,9).Value = 0.07303027402976368
$ws.Cells.Item(8,10).Value = 0.07303027402976367
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 26.04517333333333
$ws.Cells.Item(8,14).Value = 78.13552
$ws.Cells.Item(8,15).Value = 0.9210237118384171
$ws.Cells.Item(8,16).Value = 0.921023711838417
$ws.Cells.Item(8,17).Value = 36.88163233109333
$ws.Cells.Item(8,18).Value = 331.93469097984
$ws.Cells.Item(8,19).Value = 0.0672626140634697
$ws.Cells.Item(8,20).Value = 0.06726261406346969
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Mdk"
$ws.Cells.Item(9,3).Value = "Itga4"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.416064
$ws.Cells.Item(9,8).Value = 4.248192
$ws.Cells.Item(9,9).Value = 0.07303027402976368
$ws.Cells.Item(9,10).Value = 0.07303027402976367
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.3302223333333333
$ws.Cells.Item(9,14).Value = 0.990667
$ws.Cells.Item(9,15).Value = 0.01167750336256582
$ws.Cells.Item(9,16).Value = 0.01167750336256582
$ws.Cells.Item(9,17).Value = 0.4676159582293333
$ws.Cells.Item(9,18).Value = 4.208543624063999
$ws.Cells.Item(9,19).Value = 0.0008528112705516689
$ws.Cells.Item(9,20).Value = 0.0008528112705516688
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Mdk"
$ws.Cells.Item(10,3).Value = "Itga4"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.416064
$ws.Cells.Item(10,8).Value = 4.248192
$ws.Cells.Item(10,9).Value = 0.07303027402976368
$ws.Cells.Item(10,10).Value = 0.07303027402976367
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.903109
$ws.Cells.Item(10,14).Value = 5.709327
$ws.Cells.Item(10,15).Value = 0.06729878479901708
$ws.Cells.Item(10,16).Value = 0.06729878479901708
$ws.Cells.Item(10,17).Value = 2.694924142975999
$ws.Cells.Item(10,18).Value = 24.254317286784
$ws.Cells.Item(10,19).Value = 0.004914848695742311
$ws.Cells.Item(10,20).Value = 0.00491484869574231
